$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1275426186571785
$ws.Range("C2").Value = 0.9905766866993081
$ws.Range("D2").Value = 0.2759961217547767
$ws.Range("G2").Value = 0.1395347341502202
$ws.Range("H2").Value = 0.991
